$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 6.250754832778648
$ws.Cells.Item(2, 4).Value = 10.05852583752363
$ws.Cells.Item(2, 5).Value = 14.02413460890021
$ws.Cells.Item(2, 6).Value = 28.52769008755491
$ws.Cells.Item(2, 7).Value = 27.24374681691054
$ws.Cells.Item(2, 8).Value = 13.89131497501569
$ws.Cells.Item(2, 9).Value = 18.73486035391902
$ws.Cells.Item(2, 10).Value = 9.751380133820122
$ws.Cells.Item(2, 11).Value = 10.59955265147022
$ws.Cells.Item(2, 14).Value = 17.58772226885498
$ws.Cells.Item(2, 15).Value = 20.94614093855296
$ws.Cells.Item(3, 2).Value = 6.120268030835176
$ws.Cells.Item(3, 4).Value = 10.00399083386647
$ws.Cells.Item(3, 5).Value = 13.96063028545999
$ws.Cells.Item(3, 6).Value = 28.53722952529662
$ws.Cells.Item(3, 7).Value = 27.24109689909521
$ws.Cells.Item(3, 8).Value = 13.93004950772598
$ws.Cells.Item(3, 9).Value = 18.82866093917338
$ws.Cells.Item(3, 10).Value = 9.755847742273918
$ws.Cells.Item(3, 11).Value = 10.23137723200765
$ws.Cells.Item(3, 14).Value = 17.62858252090579
$ws.Cells.Item(3, 15).Value = 20.99749849808124
$ws.Cells.Item(4, 2).Value = 6.039598074390607
$ws.Cells.Item(4, 4).Value = 9.972168438534663
$ws.Cells.Item(4, 5).Value = 13.92443938601826
$ws.Cells.Item(4, 6).Value = 28.55012842108919
$ws.Cells.Item(4, 7).Value = 27.2483599233453
$ws.Cells.Item(4, 8).Value = 13.95604615145696
$ws.Cells.Item(4, 9).Value = 18.8892067014032
$ws.Cells.Item(4, 10).Value = 9.760132222165671
$ws.Cells.Item(4, 11).Value = 9.999126571361145
$ws.Cells.Item(4, 14).Value = 17.65527446811231
$ws.Cells.Item(4, 15).Value = 21.03357557029779
$ws.Cells.Item(5, 2).Value = 6.00663177690965
$ws.Cells.Item(5, 4).Value = 9.95962911411671
$ws.Cells.Item(5, 5).Value = 13.91040755653958
$ws.Cells.Item(5, 6).Value = 28.55715537053434
$ws.Cells.Item(5, 7).Value = 27.25355315505774
$ws.Cells.Item(5, 8).Value = 13.96719641089071
$ws.Cells.Item(5, 9).Value = 18.91462390546648
$ws.Cells.Item(5, 10).Value = 9.762266282348532
$ws.Cells.Item(5, 11).Value = 9.903062365789204
$ws.Cells.Item(5, 14).Value = 17.66655568018015
$ws.Cells.Item(5, 15).Value = 21.04941756555679
$ws.Cells.Item(6, 2).Value = 6.001153519477851
$ws.Cells.Item(6, 4).Value = 9.957573149571841
$ws.Cells.Item(6, 5).Value = 13.90812117177757
$ws.Cells.Item(6, 6).Value = 28.55842910355377
$ws.Cells.Item(6, 7).Value = 27.25455028963487
$ws.Cells.Item(6, 8).Value = 13.96908150350958
$ws.Cells.Item(6, 9).Value = 18.91888942807226
$ws.Cells.Item(6, 10).Value = 9.762644094448277
$ws.Cells.Item(6, 11).Value = 9.887029438867307
$ws.Cells.Item(6, 14).Value = 17.66845334253508
$ws.Cells.Item(6, 15).Value = 21.05211692093618
$ws.Cells.Item(7, 2).Value = 6.039153792886577
$ws.Cells.Item(7, 4).Value = 9.971997580089969
$ws.Cells.Item(7, 5).Value = 13.9242472327154
$ws.Cells.Item(7, 6).Value = 28.55021602103448
$ws.Cells.Item(7, 7).Value = 27.2484209218073
$ws.Cells.Item(7, 8).Value = 13.95619427510546
$ws.Cells.Item(7, 9).Value = 18.88954647026602
$ws.Cells.Item(7, 10).Value = 9.760159430850839
$ws.Cells.Item(7, 11).Value = 9.997836576580754
$ws.Cells.Item(7, 14).Value = 17.65542497353567
$ws.Cells.Item(7, 15).Value = 21.03378460698016
$ws.Cells.Item(8, 2).Value = 6.205904263585651
$ws.Cells.Item(8, 4).Value = 10.0393836646075
$ws.Cells.Item(8, 5).Value = 14.00166439989319
$ws.Cells.Item(8, 6).Value = 28.52951799780494
$ws.Cells.Item(8, 7).Value = 27.24098734877045
$ws.Cells.Item(8, 8).Value = 13.90421112884769
$ws.Cells.Item(8, 9).Value = 18.76659138522817
$ws.Cells.Item(8, 10).Value = 9.752600999119117
$ws.Cells.Item(8, 11).Value = 10.47396809992493
$ws.Cells.Item(8, 14).Value = 17.60147853410781
$ws.Cells.Item(8, 15).Value = 20.96290481519858
$ws.Cells.Item(9, 2).Value = 6.526601702435879
$ws.Cells.Item(9, 4).Value = 10.18422725342834
$ws.Cells.Item(9, 5).Value = 14.17513812730987
$ws.Cells.Item(9, 6).Value = 28.54476778076327
$ws.Cells.Item(9, 7).Value = 27.29695386784425
$ws.Cells.Item(9, 8).Value = 13.81984315951828
$ws.Cells.Item(9, 9).Value = 18.54880364402228
$ws.Cells.Item(9, 10).Value = 9.749982484382071
$ws.Cells.Item(9, 11).Value = 11.35321858162605
$ws.Cells.Item(9, 14).Value = 17.50837796979064
$ws.Cells.Item(9, 15).Value = 20.86004861193565
$ws.Cells.Item(10, 2).Value = 6.755923887157802
$ws.Cells.Item(10, 4).Value = 10.29769471841487
$ws.Cells.Item(10, 5).Value = 14.314977581435
$ws.Cells.Item(10, 6).Value = 28.58990818454216
$ws.Cells.Item(10, 7).Value = 27.3809616442963
$ws.Cells.Item(10, 8).Value = 13.76857979185798
$ws.Cells.Item(10, 9).Value = 18.402882032592
$ws.Cells.Item(10, 10).Value = 9.755457249314919
$ws.Cells.Item(10, 11).Value = 11.95975166662182
$ws.Cells.Item(10, 14).Value = 17.44766327041653
$ws.Cells.Item(10, 15).Value = 20.80662333608554
$ws.Cells.Item(11, 2).Value = 6.858405387042882
$ws.Cells.Item(11, 4).Value = 10.35068841038867
$ws.Cells.Item(11, 5).Value = 14.38109406395222
$ws.Cells.Item(11, 6).Value = 28.61776874563463
$ws.Cells.Item(11, 7).Value = 27.42842545591959
$ws.Cells.Item(11, 8).Value = 13.74758875446958
$ws.Cells.Item(11, 9).Value = 18.33952893934556
$ws.Cells.Item(11, 10).Value = 9.759543145123976
$ws.Cells.Item(11, 11).Value = 12.22602323619273
$ws.Cells.Item(11, 14).Value = 17.42170142285962
$ws.Cells.Item(11, 15).Value = 20.78714648148636
$ws.Cells.Item(12, 2).Value = 6.896911600252197
$ws.Cells.Item(12, 4).Value = 10.37093995602221
$ws.Cells.Item(12, 5).Value = 14.40647347748927
$ws.Cells.Item(12, 6).Value = 28.62936671271197
$ws.Cells.Item(12, 7).Value = 27.44771974175245
$ws.Cells.Item(12, 8).Value = 13.73997503108854
$ws.Cells.Item(12, 9).Value = 18.31597198826435
$ws.Cells.Item(12, 10).Value = 9.761318588285135
$ws.Cells.Item(12, 11).Value = 12.32539008672472
$ws.Cells.Item(12, 14).Value = 17.412107952898
$ws.Cells.Item(12, 15).Value = 20.78046633869772
$ws.Cells.Item(13, 2).Value = 6.888632561109952
$ws.Cells.Item(13, 4).Value = 10.36657046641776
$ws.Cells.Item(13, 5).Value = 14.40099261498516
$ws.Cells.Item(13, 6).Value = 28.62682238501338
$ws.Cells.Item(13, 7).Value = 27.44350579714557
$ws.Cells.Item(13, 8).Value = 13.74159987438458
$ws.Cells.Item(13, 9).Value = 18.32102614873453
$ws.Cells.Item(13, 10).Value = 9.760926086176855
$ws.Cells.Item(13, 11).Value = 12.30405592287298
$ws.Cells.Item(13, 14).Value = 17.4141635157452
$ws.Cells.Item(13, 15).Value = 20.78187408509833
$ws.Cells.Item(14, 2).Value = 6.861579557908323
$ws.Cells.Item(14, 4).Value = 10.3523509110886
$ws.Cells.Item(14, 5).Value = 14.38317528018267
$ws.Cells.Item(14, 6).Value = 28.61870195481658
$ws.Cells.Item(14, 7).Value = 27.42998639763524
$ws.Cells.Item(14, 8).Value = 13.74695565119203
$ws.Cells.Item(14, 9).Value = 18.3375822175044
$ws.Cells.Item(14, 10).Value = 9.75968464579435
$ws.Cells.Item(14, 11).Value = 12.2342279666643
$ws.Cells.Item(14, 14).Value = 17.42090740169362
$ws.Cells.Item(14, 15).Value = 20.7865829605677
$ws.Cells.Item(15, 2).Value = 6.844968518631401
$ws.Cells.Item(15, 4).Value = 10.34366456281584
$ws.Cells.Item(15, 5).Value = 14.37230573108773
$ws.Cells.Item(15, 6).Value = 28.61386422089423
$ws.Cells.Item(15, 7).Value = 27.42187707322209
$ws.Cells.Item(15, 8).Value = 13.75027986895763
$ws.Cells.Item(15, 9).Value = 18.34777968651382
$ws.Cells.Item(15, 10).Value = 9.758953908212243
$ws.Cells.Item(15, 11).Value = 12.19126347615358
$ws.Cells.Item(15, 14).Value = 17.42506916577605
$ws.Cells.Item(15, 15).Value = 20.78955786694498
$ws.Cells.Item(16, 2).Value = 6.749186272777768
$ws.Cells.Item(16, 4).Value = 10.29425797084205
$ws.Cells.Item(16, 5).Value = 14.31070560614238
$ws.Cells.Item(16, 6).Value = 28.588234387325
$ws.Cells.Item(16, 7).Value = 27.37804513808162
$ws.Cells.Item(16, 8).Value = 13.76999848969889
$ws.Cells.Item(16, 9).Value = 18.40708306689511
$ws.Cells.Item(16, 10).Value = 9.755222229238864
$ws.Cells.Item(16, 11).Value = 11.94214927799567
$ws.Cells.Item(16, 14).Value = 17.44939323790573
$ws.Cells.Item(16, 15).Value = 20.80799343127824
$ws.Cells.Item(17, 2).Value = 6.689929052136143
$ws.Cells.Item(17, 4).Value = 10.26429139268996
$ws.Cells.Item(17, 5).Value = 14.27354479409894
$ws.Cells.Item(17, 6).Value = 28.57438397133256
$ws.Cells.Item(17, 7).Value = 27.35351852386504
$ws.Cells.Item(17, 8).Value = 13.78269189015252
$ws.Cells.Item(17, 9).Value = 18.44423783754045
$ws.Cells.Item(17, 10).Value = 9.753340763898089
$ws.Cells.Item(17, 11).Value = 11.78679746480177
$ws.Cells.Item(17, 14).Value = 17.46473935374056
$ws.Cells.Item(17, 15).Value = 20.8205402350123
$ws.Cells.Item(18, 2).Value = 6.655674890364847
$ws.Cells.Item(18, 4).Value = 10.24718585644695
$ws.Cells.Item(18, 5).Value = 14.25240753459406
$ws.Cells.Item(18, 6).Value = 28.56710754291485
$ws.Cells.Item(18, 7).Value = 27.3402826288383
$ws.Cells.Item(18, 8).Value = 13.79021199559552
$ws.Cells.Item(18, 9).Value = 18.46589333931441
$ws.Cells.Item(18, 10).Value = 9.752408842235216
$ws.Cells.Item(18, 11).Value = 11.69653952452035
$ws.Cells.Item(18, 14).Value = 17.47372208966382
$ws.Cells.Item(18, 15).Value = 20.82821106491548
$ws.Cells.Item(19, 2).Value = 6.644048821096139
$ws.Cells.Item(19, 4).Value = 10.24141704157062
$ws.Cells.Item(19, 5).Value = 14.24529198689387
$ws.Cells.Item(19, 6).Value = 28.56476252671593
$ws.Cells.Item(19, 7).Value = 27.33595105012228
$ws.Cells.Item(19, 8).Value = 13.79279581401494
$ws.Cells.Item(19, 9).Value = 18.47327453596694
$ws.Cells.Item(19, 10).Value = 9.752119149494112
$ws.Cells.Item(19, 11).Value = 11.66582705654214
$ws.Cells.Item(19, 14).Value = 17.47679031216529
$ws.Cells.Item(19, 15).Value = 20.8308862495822
$ws.Cells.Item(20, 2).Value = 6.696255055354015
$ws.Cells.Item(20, 4).Value = 10.26746798574967
$ws.Cells.Item(20, 5).Value = 14.27747625679189
$ws.Cells.Item(20, 6).Value = 28.57578700073116
$ws.Cells.Item(20, 7).Value = 27.35603932111771
$ws.Cells.Item(20, 8).Value = 13.78131796816164
$ws.Cells.Item(20, 9).Value = 18.44025315950375
$ws.Cells.Item(20, 10).Value = 9.753525507415551
$ws.Cells.Item(20, 11).Value = 11.80342905242486
$ws.Cells.Item(20, 14).Value = 17.46308958697323
$ws.Cells.Item(20, 15).Value = 20.81915758418374
$ws.Cells.Item(21, 2).Value = 6.869534140368243
$ws.Cells.Item(21, 4).Value = 10.35652265883794
$ws.Cells.Item(21, 5).Value = 14.38839950826515
$ws.Cells.Item(21, 6).Value = 28.62105873638577
$ws.Cells.Item(21, 7).Value = 27.43392161034765
$ws.Cells.Item(21, 8).Value = 13.74537343275772
$ws.Cells.Item(21, 9).Value = 18.33270754858053
$ws.Cells.Item(21, 10).Value = 9.76004310389086
$ws.Cells.Item(21, 11).Value = 12.25477844532858
$ws.Cells.Item(21, 14).Value = 17.41892011233341
$ws.Cells.Item(21, 15).Value = 20.78518097014578
$ws.Cells.Item(22, 2).Value = 6.981010401109533
$ws.Cells.Item(22, 4).Value = 10.41579075761339
$ws.Cells.Item(22, 5).Value = 14.46288237927241
$ws.Cells.Item(22, 6).Value = 28.6567502713252
$ws.Cells.Item(22, 7).Value = 27.49251499908464
$ws.Cells.Item(22, 8).Value = 13.72383506611147
$ws.Cells.Item(22, 9).Value = 18.26494614558955
$ws.Cells.Item(22, 10).Value = 9.765632215350108
$ws.Cells.Item(22, 11).Value = 12.54119622473598
$ws.Cells.Item(22, 14).Value = 17.39143812465204
$ws.Cells.Item(22, 15).Value = 20.76702854300119
$ws.Cells.Item(23, 2).Value = 6.921687097249375
$ws.Cells.Item(23, 4).Value = 10.38406543873731
$ws.Cells.Item(23, 5).Value = 14.42295344094992
$ws.Cells.Item(23, 6).Value = 28.63714468709453
$ws.Cells.Item(23, 7).Value = 27.46054224539624
$ws.Cells.Item(23, 8).Value = 13.73515168597631
$ws.Cells.Item(23, 9).Value = 18.30088116875877
$ws.Cells.Item(23, 10).Value = 9.762527988609655
$ws.Cells.Item(23, 11).Value = 12.38913659509639
$ws.Cells.Item(23, 14).Value = 17.40597922953008
$ws.Cells.Item(23, 15).Value = 20.77634561196029
$ws.Cells.Item(24, 2).Value = 6.693395647052417
$ws.Cells.Item(24, 4).Value = 10.2660314649369
$ws.Cells.Item(24, 5).Value = 14.27569813417036
$ws.Cells.Item(24, 6).Value = 28.57515055275758
$ws.Cells.Item(24, 7).Value = 27.35489697427466
$ws.Cells.Item(24, 8).Value = 13.78193842522934
$ws.Cells.Item(24, 9).Value = 18.44205371425387
$ws.Cells.Item(24, 10).Value = 9.753441518303301
$ws.Cells.Item(24, 11).Value = 11.79591284635015
$ws.Cells.Item(24, 14).Value = 17.4638349479366
$ws.Cells.Item(24, 15).Value = 20.81978125561513
$ws.Cells.Item(25, 2).Value = 6.440780468848311
$ws.Cells.Item(25, 4).Value = 10.14375307532752
$ws.Cells.Item(25, 5).Value = 14.12597108586449
$ws.Cells.Item(25, 6).Value = 28.53467098253527
$ws.Cells.Item(25, 7).Value = 27.27426532178796
$ws.Cells.Item(25, 8).Value = 13.84078473021948
$ws.Cells.Item(25, 9).Value = 18.60523757816467
$ws.Cells.Item(25, 10).Value = 9.749388117305715
$ws.Cells.Item(25, 11).Value = 11.1218796689802
$ws.Cells.Item(25, 14).Value = 17.53221087419121
$ws.Cells.Item(25, 15).Value = 20.88399258001193
